$wb = $excel.ActiveWorkbook

# Sheet 1 (ALC), row 9
$ws = $wb.Worksheets.Item(1)
$ws.Range("H9").Value = 68.416664
$ws.Range("I9").Value = 22
$ws.Range("K9").Value = 22
$ws.Range("M9").Value = 147

# Sheet 1 (ALC), row 51
$ws = $wb.Worksheets.Item(1)
$ws.Range("H51").Value = 47475
$ws.Range("I51").Value = 49950
$ws.Range("J51").Value = 45000
$ws.Range("K51").Value = 49950
$ws.Range("L51").Value = 45000
$ws.Range("M51").Value = -49466
$ws.Range("N51").Value = -45968

# Sheet 1 (ALC), row 64
$ws = $wb.Worksheets.Item(1)
$ws.Range("M64").ClearContents()
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0

# Sheet 1 (ALC), row 67
$ws = $wb.Worksheets.Item(1)
$ws.Range("M67").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0

# Sheet 1 (ALC), row 92
$ws = $wb.Worksheets.Item(1)
$ws.Range("H92").Value = 117.46154
$ws.Range("I92").Value = 92
$ws.Range("J92").Value = 257.5
$ws.Range("K92").Value = 92
$ws.Range("L92").Value = 257.5
$ws.Range("M92").Value = 1156
$ws.Range("N92").Value = -2753.5

# Sheet 1 (ALC), row 99
$ws = $wb.Worksheets.Item(1)
$ws.Range("H99").Value = 3594.2144
$ws.Range("I99").Value = 475
$ws.Range("K99").Value = 1425
$ws.Range("M99").Value = 73

# Sheet 1 (ALC), row 125
$ws = $wb.Worksheets.Item(1)
$ws.Range("H125").Value = 1680.4
$ws.Range("I125").Value = 1875.5
$ws.Range("J125").Value = 900
$ws.Range("K125").Value = 16879.5
$ws.Range("L125").Value = 8100
$ws.Range("M125").Value = -14419.5
$ws.Range("N125").Value = -13020

# Sheet 1 (ALC), row 129
$ws = $wb.Worksheets.Item(1)
$ws.Range("H129").Value = 1232
$ws.Range("I129").Value = 1232
$ws.Range("K129").Value = 3696
$ws.Range("M129").Value = 1304

# Sheet 1 (ALC), row 132
$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value = 24809.8
$ws.Range("I132").Value = 26024.75
$ws.Range("K132").Value = 78074.25
$ws.Range("M132").Value = -75544.25

# Sheet 2 (ARM), row 32
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 2738.7778
$ws.Range("I32").Value = 2605.7942
$ws.Range("K32").Value = 2605.7942
$ws.Range("M32").Value = -2318.7942

# Sheet 2 (ARM), row 61
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 4587.136
$ws.Range("I61").Value = 3442.3
$ws.Range("K61").Value = 3442.3
$ws.Range("M61").Value = -3230.3

# Sheet 2 (ARM), row 102
$ws = $wb.Worksheets.Item(2)
$ws.Range("H102").Value = 4667.8887
$ws.Range("I102").Value = 3794.4285
$ws.Range("J102").Value = 7725
$ws.Range("K102").Value = 3794.4285
$ws.Range("L102").Value = 7725
$ws.Range("M102").Value = -2172.4285
$ws.Range("N102").Value = -10969

# Sheet 2 (ARM), row 136
$ws = $wb.Worksheets.Item(2)
$ws.Range("H136").Value = 4587.136
$ws.Range("I136").Value = 3442.3
$ws.Range("K136").Value = 10326.9
$ws.Range("M136").Value = -7776.900000000001

# Sheet 3 (BSM), row 105
$ws = $wb.Worksheets.Item(3)
$ws.Range("N105").ClearContents()
$ws.Range("H105").Value = 1392.2
$ws.Range("I105").Value = 1392.2
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1392.2
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 354.8

# Sheet 4 (CRP), row 22
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 1066.3334
$ws.Range("I22").Value = 1166
$ws.Range("J22").Value = 966.6667
$ws.Range("K22").Value = 1166
$ws.Range("L22").Value = 966.6667
$ws.Range("M22").Value = -816
$ws.Range("N22").Value = -1666.6667

# Sheet 4 (CRP), row 31
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 6946.2104
$ws.Range("I31").Value = 3392.6
$ws.Range("J31").Value = 8215.357
$ws.Range("K31").Value = 3392.6
$ws.Range("L31").Value = 8215.357
$ws.Range("M31").Value = -3097.6
$ws.Range("N31").Value = -8805.357

# Sheet 4 (CRP), row 34
$ws = $wb.Worksheets.Item(4)
$ws.Range("H34").Value = 6946.2104
$ws.Range("I34").Value = 3392.6
$ws.Range("J34").Value = 8215.357
$ws.Range("K34").Value = 3392.6
$ws.Range("L34").Value = 8215.357
$ws.Range("M34").Value = -3190.6
$ws.Range("N34").Value = -8619.357

# Sheet 4 (CRP), row 39
$ws = $wb.Worksheets.Item(4)
$ws.Range("M39").ClearContents()
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0

# Sheet 4 (CRP), row 49
$ws = $wb.Worksheets.Item(4)
$ws.Range("M49").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0

# Sheet 4 (CRP), row 132
$ws = $wb.Worksheets.Item(4)
$ws.Range("H132").Value = 3288
$ws.Range("I132").Value = 2258.625
$ws.Range("J132").Value = 6033
$ws.Range("K132").Value = 6775.875
$ws.Range("L132").Value = 18099
$ws.Range("M132").Value = -4245.875
$ws.Range("N132").Value = -23159

# Sheet 5 (CUL), row 2
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 26.7
$ws.Range("J2").Value = 42.75
$ws.Range("L2").Value = 256.5
$ws.Range("N2").Value = -482.5

# Sheet 5 (CUL), row 23
$ws = $wb.Worksheets.Item(5)
$ws.Range("H23").Value = 299.75
$ws.Range("J23").Value = 237.25
$ws.Range("L23").Value = 711.75
$ws.Range("N23").Value = -1181.75

# Sheet 5 (CUL), row 136
$ws = $wb.Worksheets.Item(5)
$ws.Range("H136").Value = 4127.8335
$ws.Range("I136").Value = 441.75
$ws.Range("K136").Value = 1325.25
$ws.Range("M136").Value = 3774.75

# Sheet 6 (GSM), row 7
$ws = $wb.Worksheets.Item(6)
$ws.Range("H7").Value = 18002000
$ws.Range("J7").Value = 6673333.5
$ws.Range("L7").Value = 6673333.5
$ws.Range("N7").Value = -6673557.5

# Sheet 6 (GSM), row 8
$ws = $wb.Worksheets.Item(6)
$ws.Range("H8").Value = 18002000
$ws.Range("J8").Value = 6673333.5
$ws.Range("L8").Value = 6673333.5
$ws.Range("N8").Value = -6673611.5

# Sheet 6 (GSM), row 122
$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 2708
$ws.Range("I122").Value = 2416.5
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 7249.5
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -4799.5
$ws.Range("N122").Value = -13898.5

# Sheet 6 (GSM), row 135
$ws = $wb.Worksheets.Item(6)
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140

# Sheet 7 (LTW), row 42
$ws = $wb.Worksheets.Item(7)
$ws.Range("H42").Value = 39000
$ws.Range("J42").Value = 39000
$ws.Range("L42").Value = 39000
$ws.Range("N42").Value = -40126

# Sheet 7 (LTW), row 46
$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 6855.7144
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 7998
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 7998
$ws.Range("M46").Value = -3812
$ws.Range("N46").Value = -8374

# Sheet 7 (LTW), row 49
$ws = $wb.Worksheets.Item(7)
$ws.Range("H49").Value = 39000
$ws.Range("J49").Value = 39000
$ws.Range("L49").Value = 39000
$ws.Range("N49").Value = -39294

# Sheet 7 (LTW), row 100
$ws = $wb.Worksheets.Item(7)
$ws.Range("H100").Value = 4970.2856
$ws.Range("I100").Value = 2176
$ws.Range("K100").Value = 2176
$ws.Range("M100").Value = -1635

# Sheet 7 (LTW), row 135
$ws = $wb.Worksheets.Item(7)
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140

# Sheet 8 (WVR), row 58
$ws = $wb.Worksheets.Item(8)
$ws.Range("H58").Value = 25353
$ws.Range("I58").Value = 26482.5
$ws.Range("J58").Value = 23094
$ws.Range("K58").Value = 26482.5
$ws.Range("L58").Value = 23094
$ws.Range("M58").Value = -26174.5
$ws.Range("N58").Value = -23710

# Sheet 8 (WVR), row 62
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 9790
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 10987.5
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 10987.5
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -12235.5

# Sheet 8 (WVR), row 65
$ws = $wb.Worksheets.Item(8)
$ws.Range("H65").Value = 9790
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 10987.5
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 54937.5
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -61177.5

# Sheet 8 (WVR), row 98
$ws = $wb.Worksheets.Item(8)
$ws.Range("H98").Value = 30896.25
$ws.Range("J98").Value = 30896.25
$ws.Range("L98").Value = 30896.25
$ws.Range("N98").Value = -36886.25

# Sheet 8 (WVR), row 122
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 2098.0476
$ws.Range("I122").Value = 919.5
$ws.Range("K122").Value = 2758.5
$ws.Range("M122").Value = -308.5

# Sheet 8 (WVR), row 126
$ws = $wb.Worksheets.Item(8)
$ws.Range("H126").Value = 4944.76
$ws.Range("I126").Value = 2992.7856
$ws.Range("J126").Value = 7429.091
$ws.Range("K126").Value = 8978.356800000001
$ws.Range("L126").Value = 22287.273
$ws.Range("M126").Value = -6508.356800000001
$ws.Range("N126").Value = -27227.273
